$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: update Prophet/Amazon forecast columns (C:G) for rows 2-17 ---

$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 18
$ws1.Range("E2").Value = 21
$ws1.Range("F2").Value = 25
$ws1.Range("G2").Value = 30

$ws1.Range("D3").Value = 17
$ws1.Range("E3").Value = 20
$ws1.Range("F3").Value = 24
$ws1.Range("G3").Value = 30

$ws1.Range("D4").Value = 17
$ws1.Range("E4").Value = 20
$ws1.Range("F4").Value = 24
$ws1.Range("G4").Value = 30

$ws1.Range("D5").Value = 18
$ws1.Range("E5").Value = 22
$ws1.Range("F5").Value = 26
$ws1.Range("G5").Value = 32

$ws1.Range("D6").Value = 18
$ws1.Range("E6").Value = 22
$ws1.Range("F6").Value = 27
$ws1.Range("G6").Value = 34

$ws1.Range("C7").Value = 16
$ws1.Range("D7").Value = 18
$ws1.Range("E7").Value = 22
$ws1.Range("F7").Value = 27
$ws1.Range("G7").Value = 33

$ws1.Range("D8").Value = 19
$ws1.Range("E8").Value = 24
$ws1.Range("F8").Value = 28
$ws1.Range("G8").Value = 36

$ws1.Range("C9").Value = 15
$ws1.Range("D9").Value = 20
$ws1.Range("E9").Value = 24
$ws1.Range("F9").Value = 29
$ws1.Range("G9").Value = 37

$ws1.Range("D10").Value = 18
$ws1.Range("E10").Value = 22
$ws1.Range("F10").Value = 27
$ws1.Range("G10").Value = 34

$ws1.Range("D11").Value = 19
$ws1.Range("E11").Value = 23
$ws1.Range("F11").Value = 28
$ws1.Range("G11").Value = 35

$ws1.Range("D12").Value = 19
$ws1.Range("E12").Value = 24
$ws1.Range("F12").Value = 29
$ws1.Range("G12").Value = 37

$ws1.Range("D13").Value = 21
$ws1.Range("E13").Value = 26
$ws1.Range("F13").Value = 31
$ws1.Range("G13").Value = 41

$ws1.Range("D14").Value = 20
$ws1.Range("E14").Value = 24
$ws1.Range("F14").Value = 29
$ws1.Range("G14").Value = 37

$ws1.Range("D15").Value = 19
$ws1.Range("E15").Value = 23
$ws1.Range("F15").Value = 29
$ws1.Range("G15").Value = 37

$ws1.Range("D16").Value = 18
$ws1.Range("E16").Value = 22
$ws1.Range("F16").Value = 27
$ws1.Range("G16").Value = 35

$ws1.Range("D17").Value = 18
$ws1.Range("E17").Value = 22
$ws1.Range("F17").Value = 27
$ws1.Range("G17").Value = 35

# --- Summary sheet: updated aggregate metrics (stored as text, matching existing column formatting) ---
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "243"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "84"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "25"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "2"
